$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - values must stay as text, so a leading
# apostrophe forces text entry; Style is reset to Normal afterwards so the
# cell formatting matches the original (unstyled) cells.
$ws.Range("D2").Value = "'244.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'23.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.414"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05976"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.392"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8085"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.9284"
$ws.Range("D8").Style = "Normal"

# Rows 9-17: coin list shifted up by one position, with refreshed prices
# and volume labels for each row.
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1429"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "8WazirXWRX"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07436"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03384"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03042"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09336"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'3.937"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13MCDexMCB"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001605"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitForexTokenBF"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04807"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15CoinExTokenCET"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005944"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONE"

# Remaining standalone price refreshes.
$ws.Range("D18").Value = "'0.005632"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.004155"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.0009866"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Value = "'3.661"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'6.456"
$ws.Range("D23").Style = "Normal"
$ws.Range("D40").Value = "'0.03934"
$ws.Range("D40").Style = "Normal"

# Rows 41-43: another coin-list shuffle with refreshed prices/volumes.
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1075"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40BKEXTokenBKK"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002702"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003025"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42KickTokenKICK"

# Final set of standalone price refreshes (plus a label tweak on row 44).
$ws.Range("D44").Value = "'0.007339"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").Value = "'0.00005124"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.0005804"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.8555"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.002265"
$ws.Range("D49").Style = "Normal"
